$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 33, shifting existing rows 33+ down
$ws.Rows.Item(33).Insert()

# Populate the new row 33 with the new LeetCode problem entry
$ws.Cells.Item(33, 1).Value = "Binary Tree"
$ws.Cells.Item(33, 2).Value = 226
$ws.Cells.Item(33, 3).Value = "226-Invert Binary Tree"
$ws.Cells.Item(33, 4).Value = "Easy"
$ws.Cells.Item(33, 5).Value = "BFS Iterative que , swap current nodes"
$ws.Cells.Item(33, 6).Value = "O(n), O(n) time"
$ws.Cells.Item(33, 7).Value = "O(n), O(n) time"
$ws.Cells.Item(33, 13).Value = "15 minutes"

# The worksheet's Table1 (ListObject) needs to grow to cover the newly
# inserted row, since it does not automatically expand on a plain row insert.
$table = $ws.ListObjects.Item("Table1")
$newRange = $ws.Range("A2:X78")
$table.Resize($newRange)

# Leave the selection where the author ended up after the edit
[void]$ws.Range("F36").Select()
